$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Allocation Sheet")

# --- DSR Name update (M3) ---
$ws.Range("M3").Value = "Md Robiul Islam"

# --- Insert a new model row (i68 / 5550 / 5990) above row 10 in the I:K mini table,
#     shifting the existing I10:K32 block down one row into I11:K33 ---
$src = $ws.Range("I10:K32")
$dst = $ws.Range("I11:K33")
$dst.Value = $src.Value2

$ws.Range("I10").Value = "i68"
$ws.Range("J10").Value = 5550
$ws.Range("K10").Value = 5990

# --- Fill in the "Total Qnt" row values that were previously blank ---
$ws.Range("B31").Value = 970
$ws.Range("C31").Value = 1050

# --- Scroll position shifted down one row in the saved view ---
$ws.Application.ActiveWindow.ScrollRow = 37
